# Apply updated dSF (column F) values for the davies_zach 2022 dataset.
# These reflect a data "repull" where the dSF values change for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -1
    4  = 1
    6  = -4
    12 = 2
    13 = -1
    15 = -3
    17 = -2
    18 = 1
    19 = -3
    21 = -3
    22 = 1
    26 = 2
    27 = -3
    28 = -1
    29 = 3
    30 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
